# Generate Report for Archive
#
# Every row's localization status moves on from "Ready for handoff" to
# "In Translation": the Overview sheet tracks this per-language in columns
# E (zh-cn) and F (de-de), while the per-language detail sheets (zh-cn,
# de-de) track it in their "Status" column (C). After the text is refreshed
# the narrower "In Translation" label lets the Status column(s) be resized
# down to fit it.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: columns E (zh-cn) and F (de-de) hold the status per file.
$overview.Range("E2:E4").Value = "In Translation"
$overview.Range("F2:F4").Value = "In Translation"

# zh-cn / de-de detail sheets: column C ("Status") holds the same value.
$zhcn.Range("C2:C4").Value = "In Translation"
$dede.Range("C2:C4").Value = "In Translation"

# Resize the Status column(s) to fit the new, shorter text.
$overview.Columns("E").ColumnWidth = 12.5
$overview.Columns("F").ColumnWidth = 12.5
$zhcn.Columns("C").ColumnWidth = 12.5
$dede.Columns("C").ColumnWidth = 12.5
